$wb = $excel.ActiveWorkbook

# --- Trend sheet: fill in forecasted values B10:B13 ---
$trend = $wb.Worksheets.Item("Trend")
$trend.Range("B9").Copy()
$trend.Range("B10:B13").PasteSpecial(-4122)
$trend.Range("B10").Value = 412750
$trend.Range("B11").Value = 426333.33333333302
$trend.Range("B12").Value = 439916.66666666599
$trend.Range("B13").Value = 453500

# --- Forecast sheet: add a new Miles Driven input and Predicted Amount formula ---
$forecast = $wb.Worksheets.Item("Forecast")
$forecast.Range("D2").Value = 110
$forecast.Range("E2").Formula = "=FORECAST.LINEAR(D2, B2:B27, A2:A27)"

# --- Selections / active sheet to match the final saved state ---
$trend.Range("C13").Select() | Out-Null
$forecast.Range("D3").Select() | Out-Null
$forecast.Activate() | Out-Null
